$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Assigning a plain string that looks like a date (e.g. "2022-08-15") via
    # Range.Value gets auto-parsed by Excel into a real date serial number
    # (and a new number-formatted style). Routing the literal text through a
    # formula and then converting the formula to a static value keeps it as
    # plain text without ever creating/consuming a new cell style.
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

# Update existing row 3 (Augustus Kargbo's contract) with new club/dates/fee/status
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "War Men"
Set-TextValue $ws.Range("F3") "2022-08-15"
Set-TextValue $ws.Range("G3") "2025-08-14"
$ws.Range("H3").Value = "€8M"
$ws.Range("I3").Value = "Expiring Soon"

# Add new row 4 with a prefilled contract
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "Player One"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = "War Men"
Set-TextValue $ws.Range("F4") "2024-01-01"
Set-TextValue $ws.Range("G4") "2028-12-31"
$ws.Range("H4").Value = "€15M"
$ws.Range("I4").Value = "Active"
